$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.22140880831450005
$ws.Range("B1").Value = 0.22084393921515044
$ws.Range("A2").Value = -0.15871549175709365
$ws.Range("B2").Value = 0.15713324437057796
$ws.Range("A3").Value = -0.10741506937133138
$ws.Range("B3").Value = 0.10693935313082825
$ws.Range("A4").Value = -0.098939353182974088
$ws.Range("B4").Value = 0.098513352809533217
$ws.Range("A5").Value = -0.095513352838356269
$ws.Range("B5").Value = 0.094069774546508
$ws.Range("A6").Value = -0.052728166443722913
$ws.Range("B6").Value = 0.052115092620825365
$ws.Range("A7").Value = -0.042115092693183698
$ws.Range("B7").Value = 0.041960494408522564
$ws.Range("A8").Value = -0.031960494483155522
$ws.Range("B8").Value = 0.031669798797199
$ws.Range("A9").Value = -0.029669798832368866
$ws.Range("B9").Value = 0.029425553159023288
$ws.Range("A10").Value = -0.02742555319637674
$ws.Range("B10").Value = 0.027409302111820466
$ws.Range("A11").Value = -0.024409302154771773
$ws.Range("B11").Value = 0.02438178563833393
$ws.Range("A12").Value = -0.020881785684572218
$ws.Range("B12").Value = 0.020678675071708952
$ws.Range("A13").Value = -0.017178675120117504
$ws.Range("B13").Value = 0.017086143863811465
$ws.Range("A14").Value = -0.0090861439361971108
$ws.Range("B14").Value = 0.0090558956755728559
$ws.Range("A15").Value = -0.0080558957116263485
$ws.Range("B15").Value = 0.0080363567138004299
$ws.Range("A16").Value = -0.0060363567554100328
$ws.Range("B16").Value = 0.0060036763329431864
$ws.Range("A17").Value = -0.0040036763748885207
$ws.Range("B17").Value = 0.0039999999474504833
$ws.Range("A18").Value = -0.056128447490785049
$ws.Range("B18").Value = 0.056010595388396922
$ws.Range("A19").Value = -0.052010595411193794
$ws.Range("B19").Value = 0.051143439650085831
$ws.Range("A20").Value = -0.047143439679526722
$ws.Range("B20").Value = 0.046895252793278885
$ws.Range("A21").Value = -0.0040058507095768192
$ws.Range("B21").Value = 0.0039999999688768995
$ws.Range("A22").Value = -0.045718175027676722
$ws.Range("B22").Value = 0.045503035777365142
$ws.Range("A23").Value = -0.040503035812623267
$ws.Range("B23").Value = 0.040099839591465347
$ws.Range("A24").Value = -0.020099839708014322
$ws.Range("B24").Value = 0.019999999881997965
$ws.Range("A25").Value = -0.039341608133003803
$ws.Range("B25").Value = 0.039325891144509129
$ws.Range("A26").Value = -0.030325389186325324
$ws.Range("B26").Value = 0.030316736236249753
$ws.Range("A27").Value = -0.027816736270250164
$ws.Range("B27").Value = 0.027767348482161669
$ws.Range("A28").Value = -0.025767348516047228
$ws.Range("B28").Value = 0.025735915668282594
$ws.Range("A29").Value = -0.018735915730030861
$ws.Range("B29").Value = 0.018726752280025849
$ws.Range("A30").Value = 0.041273247380609224
$ws.Range("B30").Value = -0.041461357514294139
$ws.Range("A31").Value = 0.048461357455293452
$ws.Range("B31").Value = -0.048576813700377031
$ws.Range("A32").Value = -0.0040012615824913667
$ws.Range("B32").Value = 0.003999999957843059
